# Add a new column ("Toisen asteen pohjakoulutus suoritettu") just before the
# last column ("Pohjakoulutus maa (toinen aste)"), pushing the latter one
# column to the right and leaving the new column's data cell (and the old
# column's data cell) empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at V - this shifts the existing V column
# (header "Pohjakoulutus maa (toinen aste)" / data "FIN") to W, carrying its
# formatting along automatically.
$ws.Columns("V").Insert()

# New header for the inserted column.
$ws.Range("V1").Value = "Toisen asteen pohjakoulutus suoritettu"

# The data row for the new column (V2) is left blank; the data row for the
# shifted column (W2, which used to hold "FIN") is cleared too, and its
# number format switched to text ("@"), matching the other blanked-out
# columns in that row.
$ws.Range("W2").ClearContents()
$ws.Range("W2").NumberFormat = "@"

# Move the selection to the newly added header cell.
$ws.Range("V1").Select()
